# Monthly rollover update for "CASTRO ALCIVAR EDA MARIA" workbook.
# - "VENTAS POR GRUPO": zero out cells belonging to the month that rolled off.
# - "VENTA MENSUAL": shift month columns left (Sep..Dic -> Oct..Ene), new last
#   month (enero) starts at 0; refresh totals/column widths accordingly.

$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: VENTAS POR GRUPO
# ============================================================
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCellsGrupo = @(
    "K4", "D10", "M12", "E15", "L15", "M15", "M18", "D26",
    "L26", "I29", "D30", "D31", "K31", "M31", "D33", "I33",
    "D35", "L36", "M36", "D43", "O43", "K47", "I51", "M51",
    "M52", "L55", "K61", "M61"
)
foreach ($addr in $zeroCellsGrupo) {
    $wsGrupo.Range($addr).Value = 0
}

# Row 63 totals are formatted as "<count> de 61" text, also reset to 0
$wsGrupo.Range("D63").Value = "0 de 61"
$wsGrupo.Range("E63").Value = "0 de 61"
$wsGrupo.Range("I63").Value = "0 de 61"
$wsGrupo.Range("K63").Value = "0 de 61"
$wsGrupo.Range("L63").Value = "0 de 61"
$wsGrupo.Range("M63").Value = "0 de 61"
$wsGrupo.Range("O63").Value = "0 de 61"

# ============================================================
# Sheet 2: VENTA MENSUAL
# ============================================================
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Header row: month labels shift by one (septiembre..diciembre -> octubre..enero)
$wsMensual.Range("C1").Value = "octubre"
$wsMensual.Range("D1").Value = "noviembre"
$wsMensual.Range("E1").Value = "diciembre"
$wsMensual.Range("F1").Value = "enero"

# Column widths for C/D/F narrow (E is unchanged at 15)
# iron_native/Excel stores width = ColumnWidth + 5/6, so subtract 5/6 to land on
# the exact integer widths recorded in the target file (14, 15, 11).
$wsMensual.Range("C1").ColumnWidth = 14 - 5/6
$wsMensual.Range("D1").ColumnWidth = 15 - 5/6
$wsMensual.Range("F1").ColumnWidth = 11 - 5/6

# Data rows: each advisor/client row shifts C<-D<-E<-F<-(new data), with the new
# trailing month (enero) taken from the target values below.
$mensualValues = @{
    "C4" = 1405.49
    "D4" = 461.77
    "E4" = 932.11
    "F4" = 0
    "C5" = 21118.81
    "D5" = 0
    "C6" = 387.91
    "D6" = 384.3
    "E6" = 0
    "D8" = 616.84
    "E8" = 0
    "D10" = 960.96
    "E10" = 1182.72
    "F10" = 0
    "C12" = 3404.07
    "D12" = 1466.07
    "E12" = 73.88
    "F12" = 0
    "C15" = 8398.709999999999
    "D15" = 8999.559999999999
    "E15" = 7121.26
    "F15" = 0
    "D16" = 318.84
    "E16" = 0
    "E18" = 3055.1
    "F18" = 0
    "C26" = 0
    "D26" = 3334.1
    "E26" = 5104.25
    "F26" = 0
    "C28" = 447.78
    "D28" = 1330.56
    "E28" = 0
    "C29" = 1122.59
    "D29" = 3314.82
    "E29" = 188.51
    "F29" = 0
    "C30" = 3503.69
    "D30" = 1757.95
    "E30" = 1593.98
    "F30" = 0
    "C31" = 366.83
    "D31" = 2115.57
    "E31" = 784.48
    "F31" = 0
    "C32" = 681.79
    "D32" = 0
    "C33" = 9123.92
    "D33" = 11596.22
    "E33" = 6598.56
    "F33" = 0
    "C34" = 739.2
    "D34" = 6546.97
    "E34" = 5186.46
    "D35" = 2242.43
    "E35" = 591.36
    "F35" = 0
    "C36" = 10282.96
    "D36" = 11698.32
    "E36" = 7260.09
    "F36" = 0
    "D37" = 1329.6
    "E37" = 0
    "E38" = 11545.87
    "D41" = 139.39
    "E41" = 0
    "C43" = 746.3
    "D43" = 0
    "E43" = 2701.93
    "F43" = 0
    "C44" = 9587
    "D44" = 0
    "E44" = 1994.83
    "C45" = 0
    "D45" = 4248.33
    "E45" = 0
    "C47" = 0
    "D47" = 860.8
    "E47" = 304.56
    "F47" = 0
    "D49" = 1511.66
    "E49" = 71.94
    "C50" = 879.99
    "D50" = -1206.74
    "E50" = 4401.22
    "C51" = 5784.88
    "D51" = 6643.05
    "E51" = 4861.36
    "F51" = 0
    "C52" = 6664.04
    "D52" = 4884.83
    "E52" = 5250.6
    "F52" = 0
    "C53" = 500.82
    "D53" = 388.8
    "E53" = 0
    "C55" = 4657.32
    "D55" = 0
    "E55" = 504.94
    "F55" = 0
    "C58" = 0
    "D58" = 3019.64
    "E58" = 0
    "C59" = 1035.56
    "D59" = 1101.27
    "E59" = 0
    "D60" = 8113.32
    "E60" = 0
    "C61" = 1034.72
    "D61" = 1480.85
    "E61" = 292.55
    "F61" = 0
    "C63" = 91874.38
    "D63" = 89660.08
    "E63" = 71602.56
    "F63" = 0
}
foreach ($addr in $mensualValues.Keys) {
    $wsMensual.Range($addr).Value = $mensualValues[$addr]
}
